$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp shown in A1
$ws.Range("A1").Value = "Datos actualizados a 23 de Julio de 2020 a las 05:40"

# Honduras gets fresh data and moves up (right after Barein), pushing
# Afganistan and Armenia down one row each (they keep their own prior data).
$ws.Range("A52").Value = "Honduras"
$ws.Range("B52").Value = 36102
$ws.Range("C52").Value = 757
$ws.Range("D52").Value = 4315
$ws.Range("E52").Value = 30781
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = 18
$ws.Range("H52").Value = 1006

$ws.Range("A53").Value = "Afganistan"
$ws.Range("B53").Value = 35727
$ws.Range("C53").Value = 0
$ws.Range("D53").Value = 23924
$ws.Range("E53").Value = 10613
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 1190

$ws.Range("A54").Value = "Armenia"
$ws.Range("B54").Value = 35693
$ws.Range("C54").Value = 0
$ws.Range("D54").Value = 24766
$ws.Range("E54").Value = 10249
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = 0
$ws.Range("H54").Value = 678

# Australia (row 74) data refresh
$ws.Range("D74").Value = 8656
$ws.Range("E74").Value = 4513

# Haiti (row 90) data refresh
$ws.Range("B90").Value = 7167
$ws.Range("C90").Value = 21
$ws.Range("E90").Value = 2918

# Antigua y Barbuda (row 191) data refresh
$ws.Range("D191").Value = 58
$ws.Range("E191").Value = 15

# San Vicente y las Granadinas (row 193) data refresh
$ws.Range("B193").Value = 52
$ws.Range("D193").Value = 37
$ws.Range("E193").Value = 15
